# Poll normaliser.xlsx - update booth figures and derived pref-flow medians
# for the 2026 Vic fed/state correlation model, per commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: raw booth vote counts (A:I) -------------------------------
$ws.Range("A2").Value = 31
$ws.Range("B2").Value = 33
$ws.Range("C2").Value = 15
$ws.Range("D2").Value = 10
# F2 (UAP) now has no data for this booth set -> literal #N/A error,
# matching the existing E2/G2/H2 error cells.
$ws.Range("F2").Formula = "#N/A"
$ws.Range("I2").Value = 5

# --- Row 7: preference-flow / survival medians, non-QLD branch -------
$ws.Range("C7").Formula = "=IF(`$M`$5=""qld"",0.801,0.8819)"
$ws.Range("D7").Formula = "=IF(`$M`$5=""qld"",0.338,0.255)"
$ws.Range("F7").Formula = "=IF(`$M`$5=""qld"",0.474,0.3619)"
$ws.Range("I7").Formula = "=IF(`$M`$5=""qld"",0.45,0.545)"

# --- Active cell/selection at time of save ----------------------------
$ws.Range("L17").Select()
